$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 177 (shifts existing rows 177-236 down to 178-237,
# carrying formatting/styles along as Excel normally does).
$ws.Rows.Item(177).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A177").Value = 7
$ws.Range("B177").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C177").Value = "Ñuble"
$ws.Range("D177").Value = 44841
$ws.Range("E177").Value = 16
$ws.Range("F177").Value = "Fruta"
$ws.Range("G177").Value = 100104
$ws.Range("H177").Value = "Frutos de pepita"
$ws.Range("I177").Value = 100104005
$ws.Range("J177").Value = "Pera"
$ws.Range("K177").Value = "Packham's Triumph"
$ws.Range("L177").Value = "Primera"
$ws.Range("M177").Value = 120
$ws.Range("N177").Value = 11000
$ws.Range("O177").Value = 12000
$ws.Range("P177").Value = 11500
$ws.Range("Q177").Value = "$/caja 16 kilos empedrada"
$ws.Range("R177").Value = "Provincia de Curicó"
$ws.Range("S177").Value = 719
$ws.Range("T177").Value = 16
